$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between row 16 and row 20
# so the list is sorted ascending by period (2303..2307) instead of descending.
$ws.Range("E16").Value = "2303"
$ws.Range("F16").Value = 60000
$ws.Range("E17").Value = "2304"
$ws.Range("E18").Value = "2305"
$ws.Range("E19").Value = "2306"
$ws.Range("E20").Value = "2307"
$ws.Range("F20").Value = 52000
